$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.544052880380832
$ws.Range("C2").Value = 0.109194462380934
$ws.Range("D2").Value = 0.1420870195477448
$ws.Range("F2").Value = 1.970106693555437
$ws.Range("G2").Value = 0.002514897477819356
$ws.Range("I2").Value = 1.294185951249531
$ws.Range("J2").Value = 0.2266373562450124
$ws.Range("K2").Value = 0.3130208296777539
$ws.Range("L2").Value = 0.346456203045534
$ws.Range("M2").Value = 0.2020221111795522
$ws.Range("O2").Value = 5.127613954010357
$ws.Range("B3").Value = 0.5107020152577775
$ws.Range("C3").Value = 0.1074082402434513
$ws.Range("D3").Value = 0.1406120652089413
$ws.Range("F3").Value = 1.980169137193968
$ws.Range("G3").Value = 0.002517274981234144
$ws.Range("I3").Value = 1.304014505332301
$ws.Range("J3").Value = 0.2280470300602158
$ws.Range("K3").Value = 0.2813979809856022
$ws.Range("L3").Value = 0.3443318053242095
$ws.Range("M3").Value = 0.1957320833141871
$ws.Range("O3").Value = 5.160799965435857
$ws.Range("B4").Value = 0.4903629570241606
$ws.Range("C4").Value = 0.1062998747161998
$ws.Range("D4").Value = 0.1397523330943642
$ws.Range("F4").Value = 1.987159608823063
$ws.Range("G4").Value = 0.002518813839823748
$ws.Range("I4").Value = 1.310552747517324
$ws.Range("J4").Value = 0.2289670855002797
$ws.Range("K4").Value = 0.2620010781769508
$ws.Range("L4").Value = 0.3431485506743286
$ws.Range("M4").Value = 0.1919386739211895
$ws.Range("O4").Value = 5.183332688975668
$ws.Range("B5").Value = 0.4821101143323574
$ws.Range("C5").Value = 0.1058453046469126
$ws.Range("D5").Value = 0.1394135914779113
$ws.Range("F5").Value = 1.990212732729162
$ws.Range("G5").Value = 0.002519460875702515
$ws.Range("I5").Value = 1.313343843177243
$ws.Range("J5").Value = 0.22935574515075
$ws.Range("K5").Value = 0.2541021463615039
$ws.Range("L5").Value = 0.3426969434719211
$ws.Range("M5").Value = 0.1904102539501338
$ws.Range("O5").Value = 5.193057653636089
$ws.Range("B6").Value = 0.4807418989637426
$ws.Range("C6").Value = 0.1057696489250119
$ws.Range("D6").Value = 0.1393580465987796
$ws.Range("F6").Value = 1.990732057052043
$ws.Range("G6").Value = 0.002519569521453604
$ws.Range("I6").Value = 1.313814958543762
$ws.Range("J6").Value = 0.2294211116880867
$ws.Range("K6").Value = 0.2527908801087904
$ws.Range("L6").Value = 0.3426238048039778
$ws.Range("M6").Value = 0.1901575178830264
$ws.Range("O6").Value = 5.194705267640401
$ws.Range("B7").Value = 0.4902515117327937
$ws.Range("C7").Value = 0.1062937559491459
$ws.Range("D7").Value = 0.1397477176262214
$ws.Range("F7").Value = 1.987199956187325
$ws.Range("G7").Value = 0.002518822485126426
$ws.Range("I7").Value = 1.310589876028185
$ws.Range("J7").Value = 0.2289722714793978
$ws.Range("K7").Value = 0.2618945276057758
$ws.Range("L7").Value = 0.3431423361722068
$ws.Range("M7").Value = 0.1919179903595598
$ws.Range("O7").Value = 5.183461645285433
$ws.Range("B8").Value = 0.5325251391345489
$ws.Range("C8").Value = 0.1085809953457826
$ws.Range("D8").Value = 0.1415689680749423
$ws.Range("F8").Value = 1.97340781933395
$ws.Range("G8").Value = 0.002515700866398975
$ws.Range("I8").Value = 1.297470405968681
$ws.Range("J8").Value = 0.2271121118817003
$ws.Range("K8").Value = 0.3021134948064628
$ws.Range("L8").Value = 0.3456986390755148
$ws.Range("M8").Value = 0.1998391409402949
$ws.Range("O8").Value = 5.138609207614081
$ws.Range("B9").Value = 0.6164979974180085
$ws.Range("C9").Value = 0.1129734287138646
$ws.Range("D9").Value = 0.1455020098360151
$ws.Range("F9").Value = 1.952795049022704
$ws.Range("G9").Value = 0.002510204066756931
$ws.Range("I9").Value = 1.275733266319772
$ws.Range("J9").Value = 0.2238958216218254
$ws.Range("K9").Value = 0.3811197199326841
$ws.Range("L9").Value = 0.3516680141623425
$ws.Range("M9").Value = 0.215911800475979
$ws.Range("O9").Value = 5.067744565034843
$ws.Range("B10").Value = 0.67882125468401
$ws.Range("C10").Value = 0.1161433632128706
$ws.Range("D10").Value = 0.1486089792145435
$ws.Range("F10").Value = 1.941560474248043
$ws.Range("G10").Value = 0.002506542699955731
$ws.Range("I10").Value = 1.262189133768011
$ws.Range("J10").Value = 0.2217943980504558
$ws.Range("K10").Value = 0.4392293172484187
$ws.Range("L10").Value = 0.3566313124179601
$ws.Range("M10").Value = 0.2280428013658096
$ws.Range("O10").Value = 5.026075354328384
$ws.Range("B11").Value = 0.7073049387301182
$ws.Range("C11").Value = 0.1175729242778374
$ws.Range("D11").Value = 0.1500689957414352
$ws.Range("F11").Value = 1.937295958289759
$ws.Range("G11").Value = 0.002504958155047063
$ws.Range("I11").Value = 1.256553049794377
$ws.Range("J11").Value = 0.2208948975269749
$ws.Range("K11").Value = 0.4656748833012045
$ws.Range("L11").Value = 0.3590135513199328
$ws.Range("M11").Value = 0.2336302549098406
$ws.Range("O11").Value = 5.009371220390818
$ws.Range("B12").Value = 0.7181094440250035
$ws.Range("C12").Value = 0.1181124548787977
$ws.Range("D12").Value = 0.1506285122901403
$ws.Range("F12").Value = 1.935802550333491
$ws.Range("G12").Value = 0.00250436972161793
$ws.Range("I12").Value = 1.254494235936967
$ws.Range("J12").Value = 0.2205623742736513
$ws.Range("K12").Value = 0.4756902933697802
$ws.Range("L12").Value = 0.3599334243002943
$ws.Range("M12").Value = 0.2357558650575271
$ws.Range("O12").Value = 5.003369141360906
$ws.Range("B13").Value = 0.7157816954632494
$ws.Range("C13").Value = 0.1179963382953275
$ws.Range("D13").Value = 0.1505077161394155
$ws.Range("F13").Value = 1.936118783211072
$ws.Range("G13").Value = 0.002504495936143868
$ws.Range("I13").Value = 1.254934283933657
$ws.Range("J13").Value = 0.2206336292433626
$ws.Range("K13").Value = 0.4735332579812166
$ws.Range("L13").Value = 0.3597345247295749
$ws.Range("M13").Value = 0.2352976452148212
$ws.Range("O13").Value = 5.00464741786945
$ws.Range("B14").Value = 0.7081934681035875
$ws.Range("C14").Value = 0.1176173482470659
$ws.Range("D14").Value = 0.1501148948121198
$ws.Range("F14").Value = 1.937170661410164
$ws.Range("G14").Value = 0.002504909511941392
$ws.Range("I14").Value = 1.256382158264969
$ws.Range("J14").Value = 0.2208673784495652
$ws.Range("K14").Value = 0.4664988393415115
$ws.Range("L14").Value = 0.3590888743521532
$ws.Range("M14").Value = 0.2338049353881075
$ws.Range("O14").Value = 5.008870945953504
$ws.Range("B15").Value = 0.703547828881824
$ws.Range("C15").Value = 0.1173849691316349
$ws.Range("D15").Value = 0.1498751431448682
$ws.Range("F15").Value = 1.937830780340974
$ws.Range("G15").Value = 0.002505164348058754
$ws.Range("I15").Value = 1.25727884652192
$ws.Range("J15").Value = 0.221011610633715
$ws.Range("K15").Value = 0.4621901745592254
$ws.Range("L15").Value = 0.3586957056031679
$ws.Range("M15").Value = 0.2328918741898249
$ws.Range("O15").Value = 5.011500086982721
$ws.Range("B16").Value = 0.6769623853922724
$ws.Range("C16").Value = 0.1160496854903172
$ws.Range("D16").Value = 0.1485144965090797
$ws.Range("F16").Value = 1.941856179674957
$ws.Range("G16").Value = 0.002506647879777679
$ws.Range("I16").Value = 1.26256802719589
$ws.Range("J16").Value = 0.221854316766315
$ws.Range("K16").Value = 0.4375012147283428
$ws.Range("L16").Value = 0.3564781207553267
$ws.Range("M16").Value = 0.2276790229672301
$ws.Range("O16").Value = 5.027212284067275
$ws.Range("B17").Value = 0.6606865174078678
$ws.Range("C17").Value = 0.115227325952489
$ws.Range("D17").Value = 0.1476916820169123
$ws.Range("F17").Value = 1.944542205952054
$ws.Range("G17").Value = 0.00250757869426754
$ws.Range("I17").Value = 1.265947233091236
$ws.Range("J17").Value = 0.2223857336876485
$ws.Range("K17").Value = 0.4223578084746862
$ws.Range("L17").Value = 0.3551494837539337
$ws.Range("M17").Value = 0.224498666731634
$ws.Range("O17").Value = 5.037427607227613
$ws.Range("B18").Value = 0.6513375841884965
$ws.Range("C18").Value = 0.1147531551757339
$ws.Range("D18").Value = 0.1472228149068684
$ws.Range("F18").Value = 1.946166792879566
$ws.Range("G18").Value = 0.002508121704208792
$ws.Range("I18").Value = 1.267940302686068
$ws.Range("J18").Value = 0.2226967048892252
$ws.Range("K18").Value = 0.4136488113900896
$ws.Range("L18").Value = 0.3543970040302327
$ws.Range("M18").Value = 0.2226759140998453
$ws.Range("O18").Value = 5.043515123207925
$ws.Range("B19").Value = 0.648174367078127
$ws.Range("C19").Value = 0.1145924085872778
$ws.Range("D19").Value = 0.1470648212773256
$ws.Range("F19").Value = 1.94673053723546
$ws.Range("G19").Value = 0.002508306870079274
$ws.Range("I19").Value = 1.268623615442262
$ws.Range("J19").Value = 0.2228029077034783
$ws.Range("K19").Value = 0.4107002990085675
$ws.Range("L19").Value = 0.3541442433877648
$ws.Range("M19").Value = 0.2220598833349072
$ws.Range("O19").Value = 5.045612661301192
$ws.Range("B20").Value = 0.6624178211317258
$ws.Range("C20").Value = 0.1153149889129423
$ws.Range("D20").Value = 0.1477788177437702
$ws.Range("F20").Value = 1.944248031280964
$ws.Range("G20").Value = 0.002507478818123464
$ws.Range("I20").Value = 1.265582394575951
$ws.Range("J20").Value = 0.2223286135834694
$ws.Range("K20").Value = 0.4239697412434964
$ws.Range("L20").Value = 0.3552897076145172
$ws.Range("M20").Value = 0.2248365490821413
$ws.Range("O20").Value = 5.036318236862741
$ws.Range("B21").Value = 0.7104218219477616
$ws.Range("C21").Value = 0.117728716192218
$ws.Range("D21").Value = 0.1502300962832663
$ws.Range("F21").Value = 1.936858404199903
$ws.Range("G21").Value = 0.002504787720245641
$ws.Range("I21").Value = 1.255954835419296
$ws.Range("J21").Value = 0.2207985010129558
$ws.Range("K21").Value = 0.4685649946638648
$ws.Range("L21").Value = 0.359278036172995
$ws.Range("M21").Value = 0.2342431167589325
$ws.Range("O21").Value = 5.007621618676694
$ws.Range("B22").Value = 0.7419018571852405
$ws.Range("C22").Value = 0.1192956480926597
$ws.Range("D22").Value = 0.1518708204154251
$ws.Range("F22").Value = 1.932736791262727
$ws.Range("G22").Value = 0.002503096519233107
$ws.Range("I22").Value = 1.250102407544439
$ws.Range("J22").Value = 0.2198456770383359
$ws.Range("K22").Value = 0.4977164187612004
$ws.Range("L22").Value = 0.3619881661201418
$ws.Range("M22").Value = 0.240447683509359
$ws.Range("O22").Value = 4.990751669002009
$ws.Range("B23").Value = 0.7250908440981902
$ws.Range("C23").Value = 0.11846032211362
$ws.Range("D23").Value = 0.1509916186561355
$ws.Range("F23").Value = 1.93487186472143
$ws.Range("G23").Value = 0.002503992977997313
$ws.Range("I23").Value = 1.25318574726672
$ws.Range("J23").Value = 0.2203499052076872
$ws.Range("K23").Value = 0.482157406395828
$ws.Range("L23").Value = 0.3605322853971273
$ws.Range("M23").Value = 0.2371310420327362
$ws.Range("O23").Value = 4.999583116054481
$ws.Range("B24").Value = 0.6616350728737643
$ws.Range("C24").Value = 0.1152753608015686
$ws.Range("D24").Value = 0.1477394106565413
$ws.Range("F24").Value = 1.944380777335333
$ws.Range("G24").Value = 0.002507523947741964
$ws.Range("I24").Value = 1.265747181321771
$ws.Range("J24").Value = 0.222354420597096
$ws.Range("K24").Value = 0.4232409952776379
$ws.Range("L24").Value = 0.3552262769782573
$ws.Range("M24").Value = 0.2246837747779935
$ws.Range("O24").Value = 5.03681911474149
$ws.Range("B25").Value = 0.5936688470313811
$ws.Range("C25").Value = 0.1117951593765554
$ws.Range("D25").Value = 0.1443996257435174
$ws.Range("F25").Value = 1.957683879383275
$ws.Range("G25").Value = 0.002511624603210381
$ws.Range("I25").Value = 1.281187226616602
$ws.Range("J25").Value = 0.2247198717243908
$ws.Range("K25").Value = 0.3597338296439716
$ws.Range("L25").Value = 0.3499513115011936
$ws.Range("M25").Value = 0.2115066420704004
$ws.Range("O25").Value = 5.08508811991112
